# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.300.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.853.37'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4555'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3904'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.29'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07909'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.011'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.39'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.874.40'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.908'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.157'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06646'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '85.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001028'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.94%  '
$ws.Range('E20').Value = '  -4.15%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.502'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.295.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.282'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.083.65'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.10'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.061'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.458'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '121.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9455'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09350'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.446'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.590'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.252'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06040'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02227'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.214'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.050'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.55%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5921'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1884'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.283'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5606'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.07'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.382'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.97%  '
$ws.Range('E49').Value = '  -5.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06736'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.12'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.99%  '
